# Rebuild the worker/period table (rows 16-36, columns C:E) so that it is
# grouped by worker (N. Doc / Nombre Trabajador) and, within each worker,
# sorted by Periodo Mora descending (2308 -> 2302) instead of the previous
# shuffled ordering. "Tipo Doc Trabajador" (B), "Valor Mora" (F) and
# "Salario Basico" (G) are unchanged for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$workers = @(
    @("1051419192", "CARLOS ALFONSO CASTILLO PAJARO"),
    @("1235038902", "SERGIO JOSE CABALLERO OSPINO"),
    @("1003344794", "ANDRES FELIPE BARRIOS ECHEVERRIA")
)

$periods = @("2308", "2307", "2306", "2305", "2304", "2303", "2302")

$row = 16
foreach ($worker in $workers) {
    $doc = $worker[0]
    $name = $worker[1]
    foreach ($period in $periods) {
        $ws.Cells.Item($row, 3).Value = $doc
        $ws.Cells.Item($row, 4).Value = $name
        $ws.Cells.Item($row, 5).Value = $period
        $row = $row + 1
    }
}
